$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update activation date (Ativação: 01/01/2015 -> 01/01/2021)
#
# A plain Find/Replace on this run causes the runtime to silently merge
# it with the following "Departamento: ..." run (which has no trailing
# <w:br/>) into a single <w:r>. To preserve the original run boundaries
# we instead rebuild the whole "Créditos-aula ... Departamento" line
# group (a single paragraph containing six runs separated by <w:br/>)
# by inserting a fresh copy of it - with the date corrected - right
# before the existing one, and then deleting the old copy.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Créditos-aula: 2")
$lineGroupStart = $rng.Start

$insertionPoint = $d.Range($lineGroupStart, $lineGroupStart)
$lineGroupXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListNumber"/></w:pPr><w:r><w:t>Créditos-aula: 2</w:t><w:br/></w:r><w:r><w:t>Créditos-trabalho: 0</w:t><w:br/></w:r><w:r><w:t>Carga horária: 30 h</w:t><w:br/></w:r><w:r><w:t>Semestre ideal: 2</w:t><w:br/></w:r><w:r><w:t>Ativação: 01/01/2021</w:t><w:br/></w:r><w:r><w:t>Departamento: Engenharia Química</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($lineGroupXml)

# locate the end of the newly inserted block, then remove the old one
# (which now immediately follows it, up to the original paragraph end)
$afterInsert = $d.Range($lineGroupStart, $d.Content.End)
$afterInsert.Find.Execute("Engenharia Química")
$newBlockEnd = $afterInsert.End

$oldParagraph = $d.Range($newBlockEnd, $newBlockEnd)
$oldParagraph.Expand(4)
$oldBlockEnd = $oldParagraph.End

$oldBlock = $d.Range($newBlockEnd, $oldBlockEnd - 1)
$oldBlock.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new responsible-lecturer bullet line before the existing
#    "5840560 - Marco Antonio Carvalho Pereira" one.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("5840560 - Marco Antonio Carvalho Pereira")
$newLecturerPoint = $d.Range($rng2.Start, $rng2.Start)
$newLecturerXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>11079086 - Herlandí de Souza Andrade</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newLecturerPoint.InsertXML($newLecturerXml)

# ---------------------------------------------------------------------
# 3. Update teaching method description
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras", $false, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", 2)

# ---------------------------------------------------------------------
# 4. Update evaluation criteria description
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Provas e trabalhos.", $false, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.", 2)

# ---------------------------------------------------------------------
# 5. Update recovery norm description
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Prova única com nota maior ou igual a 5,0 (cinco)", $false, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2)
